$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 951.25
$ws.Range("I19").Value = 960.4286
$ws.Range("K19").Value = 960.4286
$ws.Range("M19").Value = -785.4286
$ws.Range("H28").Value = 901.44446
$ws.Range("I28").Value = 901.8570999999999
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 901.8570999999999
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = -416.8570999999999
$ws.Range("N28").Value = -1870
$ws.Range("H70").Value = 2113.2856
$ws.Range("I70").Value = 2178
$ws.Range("J70").Value = 1951.5
$ws.Range("K70").Value = 6534
$ws.Range("L70").Value = 5854.5
$ws.Range("M70").Value = -6264
$ws.Range("N70").Value = -6394.5
$ws.Range("H73").Value = 2113.2856
$ws.Range("I73").Value = 2178
$ws.Range("J73").Value = 1951.5
$ws.Range("K73").Value = 6534
$ws.Range("L73").Value = 5854.5
$ws.Range("M73").Value = -5598
$ws.Range("N73").Value = -7726.5
$ws.Range("H100").Value = 8969
$ws.Range("I100").Value = 10132.167
$ws.Range("K100").Value = 10132.167
$ws.Range("M100").Value = -9591.166999999999
$ws.Range("H106").Value = 3217.6
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 1108.4783
$ws.Range("I107").Value = 919.7368
$ws.Range("J107").Value = 2005
$ws.Range("K107").Value = 919.7368
$ws.Range("L107").Value = 2005
$ws.Range("M107").Value = 1000.2632
$ws.Range("N107").Value = -5845
$ws.Range("H132").Value = 28274.895
$ws.Range("I132").Value = 28274.895
$ws.Range("K132").Value = 84824.685
$ws.Range("M132").Value = -82294.685

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2527.4443
$ws.Range("I2").Value = 1959.1538
$ws.Range("K2").Value = 1959.1538
$ws.Range("M2").Value = -1846.1538
$ws.Range("H32").Value = 3524971
$ws.Range("I32").Value = 3711265.5
$ws.Range("K32").Value = 3711265.5
$ws.Range("M32").Value = -3710978.5
$ws.Range("H45").Value = 3632.6155
$ws.Range("I45").Value = 2915.8572
$ws.Range("J45").Value = 4468.8335
$ws.Range("K45").Value = 2915.8572
$ws.Range("L45").Value = 4468.8335
$ws.Range("M45").Value = -2538.8572
$ws.Range("N45").Value = -5222.8335
$ws.Range("H110").Value = 2785.697
$ws.Range("I110").Value = 1706.7142
$ws.Range("K110").Value = 1706.7142
$ws.Range("M110").Value = 338.2858000000001
$ws.Range("H112").Value = 27646.5
$ws.Range("I112").Value = 19900
$ws.Range("J112").Value = 29195.8
$ws.Range("K112").Value = 19900
$ws.Range("L112").Value = 29195.8
$ws.Range("M112").Value = -18423
$ws.Range("N112").Value = -32149.8
$ws.Range("H116").Value = 2527.4443
$ws.Range("I116").Value = 1959.1538
$ws.Range("K116").Value = 1959.1538
$ws.Range("M116").Value = 334.8462
$ws.Range("H122").Value = 1663.0526
$ws.Range("I122").Value = 1186.3334
$ws.Range("K122").Value = 3559.0002
$ws.Range("M122").Value = -1109.0002
$ws.Range("H132").Value = 29414212
$ws.Range("I132").Value = 2368
$ws.Range("K132").Value = 7104
$ws.Range("M132").Value = -4574

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2527.4443
$ws.Range("I3").Value = 1959.1538
$ws.Range("K3").Value = 1959.1538
$ws.Range("M3").Value = -1845.1538
$ws.Range("H20").Value = 19603.2
$ws.Range("I20").Value = 19603.2
$ws.Range("K20").Value = 19603.2
$ws.Range("M20").Value = -19356.2
$ws.Range("H99").Value = 114703
$ws.Range("I99").Value = 8896.799999999999
$ws.Range("J99").Value = 155397.69
$ws.Range("K99").Value = 8896.799999999999
$ws.Range("L99").Value = 155397.69
$ws.Range("M99").Value = -7398.799999999999
$ws.Range("N99").Value = -158393.69
$ws.Range("H105").Value = 69107.37
$ws.Range("I105").Value = 2549.65
$ws.Range("K105").Value = 2549.65
$ws.Range("M105").Value = -802.6500000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2029.9048
$ws.Range("I105").Value = 1492.6471
$ws.Range("J105").Value = 4313.25
$ws.Range("K105").Value = 1492.6471
$ws.Range("L105").Value = 4313.25
$ws.Range("M105").Value = 254.3529000000001
$ws.Range("N105").Value = -7807.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H131").Value = 427855.25
$ws.Range("I131").Value = 126060.875
$ws.Range("J131").Value = 493108.1
$ws.Range("K131").Value = 378182.625
$ws.Range("L131").Value = 1479324.3
$ws.Range("M131").Value = -373142.625
$ws.Range("N131").Value = -1489404.3

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 4124.6
$ws.Range("I113").Value = 3968
$ws.Range("J113").Value = 4490
$ws.Range("K113").Value = 3968
$ws.Range("L113").Value = 4490
$ws.Range("M113").Value = -1798
$ws.Range("N113").Value = -8830
$ws.Range("H136").Value = 76871.75
$ws.Range("J136").Value = 76871.75
$ws.Range("L136").Value = 230615.25
$ws.Range("N136").Value = -235715.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2290.0715
$ws.Range("I46").Value = 1416.8
$ws.Range("J46").Value = 2775.2222
$ws.Range("K46").Value = 1416.8
$ws.Range("L46").Value = 2775.2222
$ws.Range("M46").Value = -1228.8
$ws.Range("N46").Value = -3151.2222
$ws.Range("H61").Value = 2079.5
$ws.Range("I61").Value = 1772.0555
$ws.Range("K61").Value = 1772.0555
$ws.Range("M61").Value = -1570.0555
$ws.Range("H100").Value = 376856.62
$ws.Range("J100").Value = 2280
$ws.Range("L100").Value = 2280
$ws.Range("N100").Value = -3362
$ws.Range("H110").Value = 30166.666
$ws.Range("J110").Value = 30166.666
$ws.Range("L110").Value = 30166.666
$ws.Range("N110").Value = -38346.666
$ws.Range("H113").Value = 2079.5
$ws.Range("I113").Value = 1772.0555
$ws.Range("K113").Value = 1772.0555
$ws.Range("M113").Value = 397.9445000000001
$ws.Range("H132").Value = 2718.0527
$ws.Range("I132").Value = 2718.0527
$ws.Range("K132").Value = 8154.158100000001
$ws.Range("M132").Value = -5624.158100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4499.3335
$ws.Range("J62").Value = 3498
$ws.Range("L62").Value = 3498
$ws.Range("N62").Value = -4746
$ws.Range("H65").Value = 4499.3335
$ws.Range("J65").Value = 3498
$ws.Range("L65").Value = 17490
$ws.Range("N65").Value = -23730
$ws.Range("H107").Value = 1563.6666
$ws.Range("I107").Value = 1070.0526
$ws.Range("K107").Value = 3210.1578
$ws.Range("M107").Value = -1290.1578
$ws.Range("H113").Value = 1541.1
$ws.Range("I113").Value = 680.9231
$ws.Range("J113").Value = 3138.5715
$ws.Range("K113").Value = 2042.7693
$ws.Range("L113").Value = 9415.7145
$ws.Range("M113").Value = 127.2307000000001
$ws.Range("N113").Value = -13755.7145
$ws.Range("H126").Value = 2205.875
$ws.Range("I126").Value = 1592.1428
$ws.Range("K126").Value = 4776.428400000001
$ws.Range("M126").Value = -2306.428400000001
$ws.Range("H132").Value = 2723.45
$ws.Range("I132").Value = 2723.45
$ws.Range("K132").Value = 8170.349999999999
$ws.Range("M132").Value = -5640.349999999999
$ws.Range("H136").Value = 903.9394
$ws.Range("I136").Value = 732.86664
$ws.Range("J136").Value = 2614.6667
$ws.Range("K136").Value = 2198.59992
$ws.Range("L136").Value = 7844.000100000001
$ws.Range("M136").Value = 351.4000800000003
$ws.Range("N136").Value = -12944.0001

Write-Output "done"